$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string comment text (added to sharedStrings.xml, referenced by row 26 col J)
$newComment = "I did only very rough sorting because it is in the cortex…. (no time to sort!)"

# Existing "highpass" text reused in column I for all new rows
$highpass = "highpass"

$rows = @(
    @{ Row=25; A=9861; B=43241; C=1; D=7; E=0; F=0.8; G=24; H=6 },
    @{ Row=26; A=9861; B=43241; C=2; D=7; E=0; F=0.8; G=24; H=6 },
    @{ Row=27; A=9861; B=43241; C=3; D=7; E=1; F=0.8; G=24; H=6 },
    @{ Row=28; A=9861; B=43241; C=4; D=7; E=0; F=0.8; G=24; H=6 },
    @{ Row=29; A=9861; B=43242; C=1; D=7; E=0; F=0.8; G=24; H=6 },
    @{ Row=30; A=9861; B=43242; C=4; D=7; E=0; F=0.8; G=24; H=6 }
)

# Copy the date cell's number format (style index 1, m/d/yy) so new date
# cells reuse the existing style instead of creating a brand new one.
$ws.Range("B24").Copy() | Out-Null

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, 2)).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $highpass
}

$excel.CutCopyMode = $false

# Add the new comment string only to row 26, column J
$ws.Cells.Item(26, 10).Value = $newComment

# Update selection to reflect the final active cell
$ws.Range("C30").Select()
